# Update apps_info sheet with latest app version info.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "ffmpeg" (row 15) was dropped from the tracked-apps list entirely.
# Deleting the whole row shifts every row below it up by one, which is
# exactly what the new layout needs (dimension goes from A1:F32 to A1:F31).
$ws.Rows("15").Delete()

# Bump the local-version column (C) for the apps whose version changed.
# Row numbers below are the POST-delete row numbers.
$ws.Range("C2").Value  = "0.30.0"      # visual_c
$ws.Range("C4").Value  = "2.4.1"       # audacity
$ws.Range("C5").Value  = "4.16.0"      # calibre
$ws.Range("C7").Value  = "8.5.2"       # crystaldiskinfo
$ws.Range("C9").Value  = "97.4.467"    # dropbox
$ws.Range("C13").Value = "11.99"       # exiftool
$ws.Range("C15").Value = "76.0.1"      # firefox (formerly row 16)
$ws.Range("C19").Value = "6.26"        # hwinfo
$ws.Range("C20").Value = "12.10.6.2"   # itunes
$ws.Range("C21").Value = "15.4.8"      # klite_codec
$ws.Range("C22").Value = "46.0.0"      # mkvtoolnix
$ws.Range("C23").Value = "25.0.8"      # obs
$ws.Range("C25").Value = "3.8.3"       # python
$ws.Range("C28").Value = "1.45"        # visual_studio_code
$ws.Range("C29").Value = "5.17.5"      # winscp
